# "Trocando Vakinha por Campanha no Desenho de Layout"
#
# Renomeia toda ocorrencia de "vakinha" para "campanha" nas tres planilhas
# (HEADER, CORPO, TRAILER) e atualiza a selecao/aba ativa do workbook para
# refletir o estado salvo pelo autor (aba TRAILER ativa, celula B3
# selecionada em CORPO).

$wb = $excel.ActiveWorkbook

# 1) Substitui "vakinha" -> "campanha" em todo o conteudo de texto das
#    celulas, em todas as planilhas (xlPart lookup, nao diferencia
#    maiusculas/minusculas) - isso cobre tanto os rotulos quanto as
#    descricoes (ex.: "Nome vakinha" -> "Nome campanha",
#    "Item vakinha" -> "Item campanha", etc.)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("vakinha", "campanha", -4142, 1, $false, $false, $false)
}

# 2) Seleciona B3 na planilha CORPO (estado salvo pelo autor).
$ws_corpo = $wb.Worksheets.Item("CORPO")
$ws_corpo.Range("B3").Select()

# 3) Ativa a planilha TRAILER, tornando-a a aba selecionada do workbook.
$ws_trailer = $wb.Worksheets.Item("TRAILER")
$ws_trailer.Activate()
